$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text-typed (it holds dotted/zero-padded
# numeric-looking strings like "0.530" or "67.976.76") rather than being
# auto-coerced to a floating point number by the Value setter.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "67.976.76"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "3.824.45"
$ws.Range("E3").Value = "  -1.82%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "600.65"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").Value = "168.71"
$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("D7").Value = "3.822.53"
$ws.Range("E7").Value = "  -1.88%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").Value = "6.52"
$ws.Range("E11").Value = "  +1.59%  "

$ws.Range("D12").Value = "0.463"
$ws.Range("E12").Value = "  +0.78%  "

$ws.Range("D13").Value = "0.0000272"
$ws.Range("E13").Value = "  +4.09%  "

$ws.Range("D14").Value = "37.13"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").Value = "4.459.33"
$ws.Range("E15").Value = "  -1.95%  "

$ws.Range("D16").Value = "3.830.20"
$ws.Range("E16").Value = "  -1.63%  "

$ws.Range("D17").Value = "18.82"
$ws.Range("E17").Value = "  +3.62%  "

$ws.Range("D18").Value = "68.032.14"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").Value = "7.40"
$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("E20").Value = "  +0.43%  "

$ws.Range("D21").Value = "10.79"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").Value = "469.78"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").Value = "0.742"
$ws.Range("E23").Value = "  +0.34%  "

$ws.Range("D24").Value = "0.0000152"
$ws.Range("E24").Value = "  -7.43%  "

$ws.Range("D25").Value = "83.71"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  +1.94%  "

$ws.Range("D27").Value = "12.26"
$ws.Range("E27").Value = "  +1.03%  "

$ws.Range("D28").Value = "10.41"
$ws.Range("E28").Value = "  +4.24%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").Value = "2.93"
$ws.Range("E30").Value = "  -1.25%  "

$ws.Range("D31").Value = "3.966.17"
$ws.Range("E31").Value = "  -1.96%  "

$ws.Range("D32").Value = "7.74"
$ws.Range("E32").Value = "  -1.54%  "

$ws.Range("D33").Value = "2.29"
$ws.Range("E33").Value = "  -1.22%  "

$ws.Range("D34").Value = "30.80"
$ws.Range("E34").Value = "  -1.73%  "

$ws.Range("D35").Value = "9.32"
$ws.Range("E35").Value = "  -1.24%  "

$ws.Range("D36").Value = "3.784.40"
$ws.Range("E36").Value = "  -2.10%  "

$ws.Range("D37").Value = "3.81"
$ws.Range("E37").Value = "  +3.17%  "

$ws.Range("E38").Value = "  +1.06%  "

$ws.Range("D39").Value = "5.99"
$ws.Range("E39").Value = "  +1.32%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.139"
$ws.Range("E40").Value = "  -0.72%  "

$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "1.01"
$ws.Range("E41").Value = "  -1.47%  "

$ws.Range("D43").Value = "0.320"
$ws.Range("E43").Value = "  +2.18%  "

$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").Value = "1.98"
$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("D46").Value = "8.78"
$ws.Range("E46").Value = "  +2.01%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "46.41"
$ws.Range("E47").Value = "  -1.53%  "

$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "406.91"
$ws.Range("E48").Value = "  -4.50%  "

$ws.Range("D49").Value = "0.000281"
$ws.Range("E49").Value = "  -8.44%  "

$ws.Range("D50").Value = "143.10"
$ws.Range("E50").Value = "  -0.21%  "

$ws.Range("D51").Value = "0.0361"
$ws.Range("E51").Value = "  +1.03%  "

# Restore the default (unstyled) cell style on the Price column so the
# text-forcing NumberFormat tweak above does not leak into the saved style.
$priceCol.Style = "Normal"
